$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 28 (2025Q2) with revised metrics
$ws.Range("C28").Value = 581
$ws.Range("D28").Value = 56
$ws.Range("E28").Value = 525
$ws.Range("F28").Value = 8.722741433021806
